# Apply "analise dos dados" corrections to the faltas_tratadas sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Turno" column (C) had the accent missing ("Manha" -> "Manhã").
$manhaRows = @(5,6,7,11,13,14,16,20,21,23,29,37,39,43,45,47,49,50,53,62,65,67,68,78,86,97,98,99,103,104,106,111,116,122,124,126,134,142,146,149,154,156,162,163,164,169,177,180,190,194,196,199,200,202,207,208,214,217)

foreach ($row in $manhaRows) {
    $ws.Range("C$row").Value = "Manhã"
}

# Rows where deeper analysis turned up more absences than originally recorded,
# which also flips the "Faltas_Suspeitas" flag to True.
$faltasUpdates = @{
    2   = 14
    9   = 11
    28  = 13
    44  = 12
    63  = 11
    98  = 12
    196 = 12
    201 = 13
}

foreach ($row in $faltasUpdates.Keys) {
    $ws.Range("D$row").Value = $faltasUpdates[$row]
    $ws.Range("H$row").Value = $true
}
